# ValidationTools.docx edit script
# 1. Insert new bullet "Checks number of records for submission" before the
#    "Checks field values against template domains where appropriate" bullet.
# 2. Move the <w:lastRenderedPageBreak/> marker from the start of the
#    "Domain Folder" paragraph to the start of the "Geodatabase" paragraph.
# 3. Split the " Optional Update Domains..." run and insert a _GoBack
#    bookmark (start/end) in the middle of it.
# 4. Remove the _GoBack bookmark that used to sit before "please include..."
#    and merge the runs back into one.

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$rNs = "http://schemas.openxmlformats.org/officeDocument/2006/relationships"
$qo = [char]0x201C   # “
$qc = [char]0x201D   # ”

# ---------------------------------------------------------------------
# 1) Add "Checks number of records for submission" bullet
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Checks field values against template domains where appropriate", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$targetPara = $rng.Paragraphs(1)
$targetIndex = $targetPara.Index
$targetPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs($targetIndex)
$newPara.Range.Text = "Checks number of records for submission"

# ---------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> from "Domain Folder" paragraph to
#    "Geodatabase" paragraph (both fully rewritten so the break lands
#    inside the same run as the leading "In the " text).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("In the " + $qo + "Geodatabase" + $qc, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$geoPara = $rng.Paragraphs(1)
$geoXml = "<w:p xmlns:w='$wNs'>" +
          "<w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""2""/></w:numPr></w:pPr>" +
          "<w:r><w:lastRenderedPageBreak/><w:t xml:space=""preserve"">In the </w:t></w:r>" +
          "<w:r><w:t>$qo" + "Geodatabase$qc</w:t></w:r>" +
          "<w:r><w:t xml:space=""preserve""> parameter, select the geodatabase </w:t></w:r>" +
          "<w:r><w:t>of data to be checked</w:t></w:r>" +
          "<w:r><w:t>.</w:t></w:r>" +
          "</w:p>"
$geoPara.Range.InsertXML($geoXml)

$rng = $d.Content
$rng.Find.Execute("In the " + $qo + "Domain Folder" + $qc, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$domPara = $rng.Paragraphs(1)
$domXml = "<w:p xmlns:w='$wNs'>" +
          "<w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""2""/></w:numPr></w:pPr>" +
          "<w:r><w:t xml:space=""preserve"">In the </w:t></w:r>" +
          "<w:r><w:t>$qo" + "Domain Folder$qc</w:t></w:r>" +
          "<w:r><w:t xml:space=""preserve""> parameter, select the $qo" + "Domains$qc folder.</w:t></w:r>" +
          "</w:p>"
$domPara.Range.InsertXML($domXml)

# ---------------------------------------------------------------------
# 3) "Optional Update Domains" paragraph: split the run and add the
#    _GoBack bookmark between " " and "Optional Update Domains..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("The script called", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$updPara = $rng.Paragraphs(1)
$updXml = "<w:p xmlns:w='$wNs' xmlns:r='$rNs'>" +
          "<w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""2""/></w:numPr></w:pPr>" +
          "<w:r><w:t>The script called $qo</w:t></w:r>" +
          "<w:r><w:t>7</w:t></w:r>" +
          "<w:r><w:t xml:space=""preserve""> </w:t></w:r>" +
          "<w:bookmarkStart w:id=""0"" w:name=""_GoBack""/>" +
          "<w:bookmarkEnd w:id=""0""/>" +
          "<w:r><w:t>Optional Update Domains$qc will sync your domains with the master copy on GitHub.</w:t></w:r>" +
          "<w:r><w:t xml:space=""preserve""> This tool requires internet access to </w:t></w:r>" +
          "<w:hyperlink r:id=""rId6"" w:history=""1""><w:r><w:t>https://raw.githubusercontent.com/kansasgis</w:t></w:r></w:hyperlink>" +
          "</w:p>"
$updPara.Range.InsertXML($updXml)

$rng = $d.Content
$rng.Find.Execute("https://raw.githubusercontent.com/kansasgis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Style = "Hyperlink"

# ---------------------------------------------------------------------
# 4) Remove the old _GoBack bookmark before "please include..." and
#    merge the surrounding runs back together.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("For issues or questions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$issuesPara = $rng.Paragraphs(1)
$issuesXml = "<w:p xmlns:w='$wNs' xmlns:r='$rNs'>" +
             "<w:r><w:t>For issues or questions, please contact Kristen Jordan</w:t></w:r>" +
             "<w:r><w:t>-Koenig</w:t></w:r>" +
             "<w:r><w:t xml:space=""preserve""> with the Kansas Data Access and Support Center. </w:t></w:r>" +
             "<w:r><w:t>Email Kristen at</w:t></w:r>" +
             "<w:r><w:t xml:space=""preserve""> </w:t></w:r>" +
             "<w:hyperlink r:id=""rId7"" w:history=""1""><w:r><w:t>Kristen@kgs.ku.edu</w:t></w:r></w:hyperlink>" +
             "<w:r><w:t xml:space=""preserve""> and</w:t></w:r>" +
             "<w:r><w:t xml:space=""preserve""> please include in the email which script you were running, any error messages, and a zipped copy of your geodatabase (change the file extension from zip to </w:t></w:r>" +
             "<w:proofErr w:type=""spellStart""/>" +
             "<w:r><w:t>piz</w:t></w:r>" +
             "<w:proofErr w:type=""spellEnd""/>" +
             "<w:r><w:t xml:space=""preserve""> so it gets through the email server).</w:t></w:r>" +
             "</w:p>"
$issuesPara.Range.InsertXML($issuesXml)

$rng = $d.Content
$rng.Find.Execute("Kristen@kgs.ku.edu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Style = "Hyperlink"

Write-Host "Edits applied"
